$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 83335410
$ws.Range("I106").Value = 90911180
$ws.Range("J106").Value = 1906
$ws.Range("K106").Value = 90911180
$ws.Range("L106").Value = 1906
$ws.Range("M106").Value = -90910549
$ws.Range("N106").Value = -3168

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 418530.47
$ws.Range("I129").Value = 500.73685
$ws.Range("J129").Value = 2007043.4
$ws.Range("K129").Value = 1502.21055
$ws.Range("L129").Value = 6021130.199999999
$ws.Range("M129").Value = 3497.78945
$ws.Range("N129").Value = -6031130.199999999

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 1853416.1
$ws.Range("I132").Value = 1297.7678
$ws.Range("J132").Value = 27783074
$ws.Range("K132").Value = 3893.3034
$ws.Range("L132").Value = 83349222
$ws.Range("M132").Value = -1363.3034
$ws.Range("N132").Value = -83354282

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1957.2727
$ws.Range("I141").Value = 1702.7778
$ws.Range("J141").Value = 3102.5
$ws.Range("K141").Value = 5108.3334
$ws.Range("L141").Value = 9307.5
$ws.Range("M141").Value = 71.66659999999956
$ws.Range("N141").Value = -19667.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 16819.455
$ws.Range("I2").Value = 18391.45
$ws.Range("J2").Value = 1099.5
$ws.Range("K2").Value = 18391.45
$ws.Range("L2").Value = 1099.5
$ws.Range("M2").Value = -18278.45
$ws.Range("N2").Value = -1325.5

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 8339891.5
$ws.Range("I32").Value = 1844233.5
$ws.Range("K32").Value = 1844233.5
$ws.Range("M32").Value = -1843946.5

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 589403.5600000001
$ws.Range("I45").Value = 715461.5
$ws.Range("J45").Value = 1133.3334
$ws.Range("K45").Value = 715461.5
$ws.Range("L45").Value = 1133.3334
$ws.Range("M45").Value = -715084.5
$ws.Range("N45").Value = -1887.3334

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 16819.455
$ws.Range("I116").Value = 18391.45
$ws.Range("J116").Value = 1099.5
$ws.Range("K116").Value = 18391.45
$ws.Range("L116").Value = 1099.5
$ws.Range("M116").Value = -16097.45
$ws.Range("N116").Value = -5687.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 16819.455
$ws.Range("I3").Value = 18391.45
$ws.Range("J3").Value = 1099.5
$ws.Range("K3").Value = 18391.45
$ws.Range("L3").Value = 1099.5
$ws.Range("M3").Value = -18277.45
$ws.Range("N3").Value = -1327.5

# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 22740554
$ws.Range("I20").Value = 35723172
$ws.Range("J20").Value = 20975
$ws.Range("K20").Value = 35723172
$ws.Range("L20").Value = 20975
$ws.Range("M20").Value = -35722925
$ws.Range("N20").Value = -21469

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1842.62
$ws.Range("I86").Value = 1887.3334
$ws.Range("J86").Value = 1248.5714
$ws.Range("K86").Value = 1887.3334
$ws.Range("L86").Value = 1248.5714
$ws.Range("M86").Value = -764.3334
$ws.Range("N86").Value = -3494.5714

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1842.62
$ws.Range("I89").Value = 1887.3334
$ws.Range("J89").Value = 1248.5714
$ws.Range("K89").Value = 9436.666999999999
$ws.Range("L89").Value = 6242.857
$ws.Range("M89").Value = -3820.666999999999
$ws.Range("N89").Value = -17474.857

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1184.174
$ws.Range("I99").Value = 933.5263
$ws.Range("J99").Value = 2374.75
$ws.Range("K99").Value = 933.5263
$ws.Range("L99").Value = 2374.75
$ws.Range("M99").Value = 564.4737
$ws.Range("N99").Value = -5370.75

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1726.8438
$ws.Range("I105").Value = 1670.0454
$ws.Range("J105").Value = 1851.8
$ws.Range("K105").Value = 1670.0454
$ws.Range("L105").Value = 1851.8
$ws.Range("M105").Value = 76.95460000000003
$ws.Range("N105").Value = -5345.8

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1000691.1
$ws.Range("I107").Value = 1429044.4
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1429044.4
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -1427124.4
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CRP")
# Row 45: A Tree Grew in Gridania / Pastoral Oak Cane
$ws.Range("H45").Value = 8500
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 12000
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = -4407
$ws.Range("N45").Value = -13186

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 10923.682
$ws.Range("I99").Value = 7017
$ws.Range("J99").Value = 16566.666
$ws.Range("K99").Value = 7017
$ws.Range("L99").Value = 16566.666
$ws.Range("M99").Value = -5519
$ws.Range("N99").Value = -19562.666

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 10923.682
$ws.Range("I126").Value = 7017
$ws.Range("J126").Value = 16566.666
$ws.Range("K126").Value = 21051
$ws.Range("L126").Value = 49699.99800000001
$ws.Range("M126").Value = -18581
$ws.Range("N126").Value = -54639.99800000001

$ws = $wb.Worksheets.Item("CUL")
# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 38097412
$ws.Range("J129").Value = 6495876
$ws.Range("L129").Value = 19487628
$ws.Range("N129").Value = -19497628

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 17893622
$ws.Range("J131").Value = 1033.2941
$ws.Range("L131").Value = 3099.8823
$ws.Range("N131").Value = -13179.8823

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 27780732
$ws.Range("I122").Value = 3969
$ws.Range("J122").Value = 55557492
$ws.Range("K122").Value = 11907
$ws.Range("L122").Value = 166672476
$ws.Range("M122").Value = -9457
$ws.Range("N122").Value = -166677376

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 7392.2
$ws.Range("I22").Value = 550
$ws.Range("K22").Value = 550
$ws.Range("M22").Value = -255

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 7392.2
$ws.Range("I27").Value = 550
$ws.Range("K27").Value = 550
$ws.Range("M27").Value = -443

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 1412.45
$ws.Range("I40").Value = 1462
$ws.Range("J40").Value = 1131.6666
$ws.Range("K40").Value = 1462
$ws.Range("L40").Value = 1131.6666
$ws.Range("M40").Value = -1326
$ws.Range("N40").Value = -1403.6666

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 16558.53
$ws.Range("I93").Value = 3400.4167
$ws.Range("J93").Value = 48138
$ws.Range("K93").Value = 3400.4167
$ws.Range("L93").Value = 48138
$ws.Range("M93").Value = -2152.4167
$ws.Range("N93").Value = -50634

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 7832.394
$ws.Range("I136").Value = 7243.5557
$ws.Range("J136").Value = 8539
$ws.Range("K136").Value = 21730.6671
$ws.Range("L136").Value = 25617
$ws.Range("M136").Value = -19180.6671
$ws.Range("N136").Value = -30717
